$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D2").Value = "275.45"
$ws.Range("E2").Value = "-1.10%"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D3").Value = "27.36"
$ws.Range("E3").Value = "1.76%"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.35%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06349"
$ws.Range("E5").Value = "-0.69%"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D6").Value = "6.955"
$ws.Range("E6").Value = "-0.67%"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("D7").Value = "1.347"
$ws.Range("E7").Value = "18.84%"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8782"
$ws.Range("E8").Value = "-1.04%"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1522"
$ws.Range("E9").Value = "2.07%"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05072"
$ws.Range("E10").Value = "-3.11%"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07564"
$ws.Range("E11").Value = "3.36%"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02965"
$ws.Range("E12").Value = "-4.86%"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09036"
$ws.Range("E13").Value = "-0.35%"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001570"
$ws.Range("E14").Value = "0.42%"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006400"
$ws.Range("E15").Value = "0.92%"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005734"
$ws.Range("E16").Value = "-4.79%"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.26%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D18").Value = "3.301"
$ws.Range("E18").Value = "-1.71%"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D19").Value = "2.285"
$ws.Range("E19").Value = "0.17%"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.21%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1352"
$ws.Range("E21").Value = "1.48%"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("D22").Value = "3.895"
$ws.Range("E22").Value = "-0.81%"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04413"
$ws.Range("E23").Value = "1.12%"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001170"
$ws.Range("E24").Value = "-0.85%"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.003866"
$ws.Range("E25").Value = "5.00%"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.25%"
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001934"
$ws.Range("E27").Value = "13.72%"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Style = "Normal"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04208"
$ws.Range("E40").Value = "3.27%"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006853"
$ws.Range("E41").Value = "2.97%"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1179"
$ws.Range("E42").Value = "0.22%"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002048"
$ws.Range("E43").Value = "-13.35%"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01149"
$ws.Range("E44").Value = "-10.80%"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005167"
$ws.Range("E45").Value = "-2.06%"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-36.78%"
$ws.Range("E46").Style = "Normal"

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D47").Value = "0.02297"
$ws.Range("E47").Value = "8.16%"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Style = "Normal"
